# Apply odds updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("J9").Value = 1.05
$ws.Range("K9").Value = 11

# Row 12
$ws.Range("G12").Value = 1.95
$ws.Range("H12").Value = 3.4
$ws.Range("I12").Value = 3.8
$ws.Range("T12").Value = 8.5
$ws.Range("AB12").Value = 13

# Row 17
$ws.Range("G17").Value = 2.38
$ws.Range("I17").Value = 2.7
$ws.Range("R17").Value = 1.44
$ws.Range("S17").Value = 2.63
$ws.Range("V17").Value = 10
$ws.Range("AJ17").Value = 21

# Row 19
$ws.Range("H19").Value = 4.15
$ws.Range("U19").Value = 25
$ws.Range("Y19").Value = 32

# Row 21
$ws.Range("N21").Value = 1.38
$ws.Range("T21").Value = 13.5
$ws.Range("U21").Value = 14
$ws.Range("W21").Value = 20
$ws.Range("Y21").Value = 16
$ws.Range("AE21").Value = 19.5
$ws.Range("AF21").Value = 26
$ws.Range("AJ21").Value = 21

# Row 25
$ws.Range("L25").Value = 1.14
$ws.Range("M25").Value = 5.5

# Row 41
$ws.Range("G41").Value = 2.7
$ws.Range("I41").Value = 2.45
$ws.Range("R41").Value = 1.67
$ws.Range("S41").Value = 2.1
$ws.Range("T41").Value = 10
$ws.Range("U41").Value = 15
$ws.Range("W41").Value = 29
$ws.Range("Z41").Value = 12
$ws.Range("AD41").Value = 151
$ws.Range("AH41").Value = 23
$ws.Range("AI41").Value = 19
$ws.Range("AJ41").Value = 26

# Row 50
$ws.Range("H50").Value = 4.15
$ws.Range("I50").Value = 4.85
$ws.Range("K50").Value = 9
$ws.Range("P50").Value = 1.27
$ws.Range("Q50").Value = 3.4
$ws.Range("V50").Value = 8
$ws.Range("X50").Value = 11.5
$ws.Range("Z50").Value = 9
$ws.Range("AE50").Value = 18
